$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- B column: mark cases 2..36 (rows 3..37) as applicable (0 -> 1) ---
for ($r = 3; $r -le 37; $r++) {
    $ws.Cells.Item($r, 2).Value = 1
}

# --- A column: highlight case blocks with a light fill + keep wrap text ---
# Rows 2..19 -> light Accent1 (blue) fill
$ws.Range("A2:A19").Interior.ThemeColor = 5
$ws.Range("A2:A19").Interior.TintAndShade = 0.79998168889431442
$ws.Range("A2:A19").WrapText = $true

# Rows 20..37 -> light Accent4 (gold) fill
$ws.Range("A20:A37").Interior.ThemeColor = 8
$ws.Range("A20:A37").Interior.TintAndShade = 0.79998168889431442
$ws.Range("A20:A37").WrapText = $true

# --- Fix the saved selection / active cell on the frozen-pane view ---
$ws.Activate()
$ws.Range("C3").Select()
